{"js": "// Apply \"Heading 1\" style to the (single, empty) paragraph and explicitly\n// clear any inherited list numbering on it, so the paragraph ends up with\n//   <w:pPr>\n//     <w:pStyle w:val=\"1\"/>         (= built-in \"Heading 1\")\n//     <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"0\"/></w:numPr>\n//   </w:pPr>\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n\n// Set the paragraph style to Heading 1 (-> w:pStyle w:val=\"1\").\nparagraph.style = \"Heading 1\";\n\n// Detach the paragraph from any numbered/bulleted list, which records an\n// explicit \"no numbering\" override (ilvl=0 / numId=0) on the paragraph.\nparagraph.detachFromList();\n\nawait context.sync();\n", "ps1": "# Apply \"Heading 1\" style to the (single, empty) paragraph and explicitly\n# remove any inherited list numbering on it, so the paragraph ends up with\n#   <w:pPr>\n#     <w:pStyle w:val=\"1\"/>         (= built-in \"Heading 1\")\n#     <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"0\"/></w:numPr>\n#   </w:pPr>\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n$p.Style = \"Heading 1\"\n$p.Range.ListFormat.RemoveNumbers()\n"}
